$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = "Figma_food_website"
$ws.Cells.Item(2, 2).Value = 77
$ws.Cells.Item(2, 3).Value = 75
$ws.Cells.Item(3, 1).Value = "Trees_website"
$ws.Cells.Item(3, 2).Value = 40
$ws.Cells.Item(3, 3).Value = 73
$ws.Cells.Item(4, 1).Value = "JoeCoffee-website-february-2022"
$ws.Cells.Item(4, 2).Value = 83
$ws.Cells.Item(4, 3).Value = 83
$ws.Cells.Item(5, 1).Value = "Hyer_website"
$ws.Cells.Item(5, 2).Value = 80
$ws.Cells.Item(5, 3).Value = 83
$ws.Cells.Item(6, 1).Value = "webflow-full"
$ws.Cells.Item(6, 2).Value = 72
$ws.Cells.Item(6, 3).Value = 80
$ws.Cells.Item(7, 1).Value = "Superlist_website"
$ws.Cells.Item(7, 2).Value = 85
$ws.Cells.Item(7, 3).Value = 85
$ws.Cells.Item(8, 1).Value = "Spotify_website"
$ws.Cells.Item(8, 2).Value = 70
$ws.Cells.Item(8, 3).Value = 75
$ws.Cells.Item(9, 1).Value = "mubasic_website"
$ws.Cells.Item(9, 2).Value = 87
$ws.Cells.Item(9, 3).Value = 83
$ws.Cells.Item(10, 1).Value = "overflow_website"
$ws.Cells.Item(10, 2).Value = 87
$ws.Cells.Item(10, 3).Value = 87
$ws.Cells.Item(11, 1).Value = "RCA_website"
$ws.Cells.Item(11, 2).Value = 60
$ws.Cells.Item(11, 3).Value = 72
$ws.Cells.Item(12, 1).Value = "Crypto_website"
$ws.Cells.Item(12, 2).Value = 90
$ws.Cells.Item(12, 3).Value = 95

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = "Figma_food_website"
$ws.Cells.Item(2, 2).Value = 75
$ws.Cells.Item(2, 3).Value = 75
$ws.Cells.Item(3, 1).Value = "Trees_website"
$ws.Cells.Item(3, 2).Value = 20
$ws.Cells.Item(3, 3).Value = 60
$ws.Cells.Item(4, 1).Value = "JoeCoffee-website-february-2022"
$ws.Cells.Item(4, 2).Value = 85
$ws.Cells.Item(4, 3).Value = 90
$ws.Cells.Item(5, 1).Value = "Hyer_website"
$ws.Cells.Item(5, 2).Value = 70
$ws.Cells.Item(5, 3).Value = 85
$ws.Cells.Item(6, 1).Value = "webflow-full"
$ws.Cells.Item(6, 2).Value = 85
$ws.Cells.Item(6, 3).Value = 80
$ws.Cells.Item(7, 1).Value = "Superlist_website"
$ws.Cells.Item(7, 2).Value = 85
$ws.Cells.Item(7, 3).Value = 85
$ws.Cells.Item(8, 1).Value = "Spotify_website"
$ws.Cells.Item(8, 2).Value = 70
$ws.Cells.Item(8, 3).Value = 70
$ws.Cells.Item(9, 1).Value = "mubasic_website"
$ws.Cells.Item(9, 2).Value = 85
$ws.Cells.Item(9, 3).Value = 85
$ws.Cells.Item(10, 1).Value = "overflow_website"
$ws.Cells.Item(10, 2).Value = 85
$ws.Cells.Item(10, 3).Value = 85
$ws.Cells.Item(11, 1).Value = "RCA_website"
$ws.Cells.Item(11, 2).Value = 70
$ws.Cells.Item(11, 3).Value = 70
$ws.Cells.Item(12, 1).Value = "Crypto_website"
$ws.Cells.Item(12, 2).Value = 85
$ws.Cells.Item(12, 3).Value = 95

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = "Figma_food_website"
$ws.Cells.Item(2, 2).Value = 85
$ws.Cells.Item(2, 3).Value = 85
$ws.Cells.Item(3, 1).Value = "Trees_website"
$ws.Cells.Item(3, 2).Value = 90
$ws.Cells.Item(3, 3).Value = 90
$ws.Cells.Item(4, 1).Value = "JoeCoffee-website-february-2022"
$ws.Cells.Item(4, 2).Value = 90
$ws.Cells.Item(4, 3).Value = 85
$ws.Cells.Item(5, 1).Value = "Hyer_website"
$ws.Cells.Item(5, 2).Value = 90
$ws.Cells.Item(5, 3).Value = 95
$ws.Cells.Item(6, 1).Value = "webflow-full"
$ws.Cells.Item(6, 2).Value = 60
$ws.Cells.Item(6, 3).Value = 90
$ws.Cells.Item(7, 1).Value = "Superlist_website"
$ws.Cells.Item(7, 2).Value = 90
$ws.Cells.Item(7, 3).Value = 90
$ws.Cells.Item(8, 1).Value = "Spotify_website"
$ws.Cells.Item(8, 2).Value = 80
$ws.Cells.Item(8, 3).Value = 95
$ws.Cells.Item(9, 1).Value = "mubasic_website"
$ws.Cells.Item(9, 2).Value = 95
$ws.Cells.Item(9, 3).Value = 90
$ws.Cells.Item(10, 1).Value = "overflow_website"
$ws.Cells.Item(10, 2).Value = 95
$ws.Cells.Item(10, 3).Value = 95
$ws.Cells.Item(11, 1).Value = "RCA_website"
$ws.Cells.Item(11, 2).Value = 60
$ws.Cells.Item(11, 3).Value = 85
$ws.Cells.Item(12, 1).Value = "Crypto_website"
$ws.Cells.Item(12, 2).Value = 100
$ws.Cells.Item(12, 3).Value = 100

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 1).Value = "Figma_food_website"
$ws.Cells.Item(2, 2).Value = 70
$ws.Cells.Item(2, 3).Value = 60
$ws.Cells.Item(3, 1).Value = "Trees_website"
$ws.Cells.Item(3, 2).Value = 30
$ws.Cells.Item(3, 3).Value = 70
$ws.Cells.Item(4, 1).Value = "JoeCoffee-website-february-2022"
$ws.Cells.Item(4, 2).Value = 75
$ws.Cells.Item(4, 3).Value = 75
$ws.Cells.Item(5, 1).Value = "Hyer_website"
$ws.Cells.Item(5, 2).Value = 80
$ws.Cells.Item(5, 3).Value = 70
$ws.Cells.Item(6, 1).Value = "webflow-full"
$ws.Cells.Item(6, 2).Value = 70
$ws.Cells.Item(6, 3).Value = 70
$ws.Cells.Item(7, 1).Value = "Superlist_website"
$ws.Cells.Item(7, 2).Value = 75
$ws.Cells.Item(7, 3).Value = 80
$ws.Cells.Item(8, 1).Value = "Spotify_website"
$ws.Cells.Item(8, 2).Value = 60
$ws.Cells.Item(8, 3).Value = 60
$ws.Cells.Item(9, 1).Value = "mubasic_website"
$ws.Cells.Item(9, 2).Value = 80
$ws.Cells.Item(9, 3).Value = 75
$ws.Cells.Item(10, 1).Value = "overflow_website"
$ws.Cells.Item(10, 2).Value = 80
$ws.Cells.Item(10, 3).Value = 80
$ws.Cells.Item(11, 1).Value = "RCA_website"
$ws.Cells.Item(11, 2).Value = 50
$ws.Cells.Item(11, 3).Value = 60
$ws.Cells.Item(12, 1).Value = "Crypto_website"
$ws.Cells.Item(12, 2).Value = 90
$ws.Cells.Item(12, 3).Value = 90
